$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff": refresh the "Latest Handoff Datetime" (column D)
# for the d004172b-70b9-449d-89e0-0618d317ef77 row (row 4) on both the zh-cn and
# de-de localization-status sheets.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-29 11:53:53"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-29 11:54:04"
